$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.550.66"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "3.613.35"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.18"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "190.43"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("D7").Value = "3.610.95"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.633"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.78%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.186"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +4.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.664"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "56.07"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -4.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000311"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +7.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.73"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.64%  "
$ws.Range("D15").Value = "4.200.75"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.98"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.57%  "
$ws.Range("D17").Value = "3.617.69"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").Value = "70.565.92"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.70"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "491.06"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "19.29"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.95"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -7.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.21"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +6.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.40"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.37%  "
$ws.Range("E27").Value = "  -4.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.09"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.47"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.40"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.61"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.28"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.119"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "66.28"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "584.03"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -8.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "39.07"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.91%  "
$ws.Range("D37").Value = "0.0₃0820"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("E38").Value = "  +0.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.400"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.27"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +20.93%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.49"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.77%  "
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.89"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.41%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "3.244.48"
$ws.Range("E43").Value = "  -1.65%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.137"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -6.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.08"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.34%  "
$ws.Range("E46").Value = "  -0.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.76"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +7.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.41"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.40%  "
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.27"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.00%  "
